# Reposition/resize the floating "Group 31" shape and switch the page
# from landscape to portrait (dimensions swapped), per the commit diff.

$d = $word.ActiveDocument

# --- Move / resize the anchored drawing (wp:anchor -> wp:posOffset / wp:extent) ---
$s = $d.Shapes.Item(1)

# Target EMUs: positionH=-481965, positionV=1226657, extent cx=6980555 cy=6113780
# (1 pt = 12700 EMU). Height (cy) is unchanged from the original 6113780 EMU = 481.4pt.
$s.Left   = -481965 / 12700          # -37.95 pt
$s.Top    = 1226657 / 12700          # 96.58716535433071 pt
$s.Width  = 549.64997                # rounds to 6980555 EMU (was 6980668)
$s.Height = 6113780 / 12700          # 481.4 pt (unchanged)

# --- Flip the page orientation: landscape (16834x11909 twips) -> portrait (11909x16834) ---
$d.PageSetup.Orientation = 0
